# Insert a new data row at row 801 (pushing the existing rows 801-842 down
# to 802-843), then populate it with the new log entry.
#   A801 = 2026/02/09   (date, stored as text like the rest of column A)
#   B801 = 月            (weekday label)
#   C801 = 6
#   D801 = 201
#
# Sheet dimension grows from A1:D842 to A1:D843 as a natural consequence of
# the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 801.. down by one to make room for the new entry.
$ws.Rows.Item(801).Insert()

# Column A in this sheet holds the date as plain text (e.g. "2026/12/29"),
# not an Excel date serial. Force text formatting before assigning so the
# "2026/02/09" string isn't auto-converted into a date value, then restore
# the default style so no stray formatting is left behind.
$ws.Range("A801").NumberFormat = "@"
$ws.Range("A801").Value = "2026/02/09"
$ws.Range("A801").Style = "Normal"

$ws.Range("B801").Value = "月"
$ws.Range("C801").Value = 6
$ws.Range("D801").Value = 201
